$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The target paragraph currently reads (three runs + a collapsed "_GoBack"
# bookmark at the very end):
#
#   "Οι γιατροί μπορούν να θεραπεύουν τους ασθενείς καταγράφοντας την
#    διάγνωσή τους" + " " + "και την προτεινόμενη θεραπεία" [_GoBack]
#
# It needs to become:
#
#   "Οι γιατροί μπορούν να θεραπεύουν τους ασθενείς καταγράφοντας την"
#   + " ημερομηνία, τη" [_GoBack] + " διάγνωσή τους" + " "
#   + "και την προτεινόμενη θεραπεία"
#
# i.e. insert " ημερομηνία, τη" right after "...καταγράφοντας την" and move
# the (hidden) "_GoBack" bookmark so that it now sits right after the newly
# typed text instead of at the end of the sentence.
# ---------------------------------------------------------------------------

# Find the exact insertion point: right after "καταγράφοντας την" and right
# before " διάγνωσή τους" (use a duplicate range so the document itself is
# not touched by the search).
$search = $d.Content.Duplicate
$search.Find.ClearFormatting()
$search.Find.Text = "καταγράφοντας την"
$search.Find.Execute() | Out-Null
$splitOffset = $search.End

# Step 1: drop a temporary bookmark exactly at that point. This alone makes
# Word split the run there (it does NOT rewrite/merge the surrounding runs),
# giving us a clean boundary to work with.
$splitPoint = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("ZZZ_TEMP_SPLIT", $splitPoint)

# Step 2: type the new text in at that same boundary.
$newText = " ημερομηνία, τη"
$insPoint = $d.Range($splitOffset, $splitOffset)
$insPoint.InsertBefore($newText)

# Step 3: the newly inserted text merges into the preceding run by default;
# force it back out into its own run (matching the target, where
# "...καταγράφοντας την" and " ημερομηνία, τη" are two separate runs) by
# round-tripping a character format on just the new span.
$newRange = $d.Range($splitOffset, $splitOffset + $newText.Length)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0

# Step 4: relocate the "_GoBack" bookmark so it now sits right after the
# newly inserted text (exactly where our temp bookmark is), instead of at
# the end of the sentence.
$tempRange = $d.Bookmarks("ZZZ_TEMP_SPLIT").Range.Duplicate
$d.Bookmarks.Add("_GoBack", $tempRange)

# Step 5: clean up the temporary bookmark.
$d.Bookmarks("ZZZ_TEMP_SPLIT").Delete()
